$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Workbook-level window/view state
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).WindowState = -4137
$excel.ActiveWindow.Left = -120
$excel.ActiveWindow.Top = -120
$excel.ActiveWindow.Width = 29040
$excel.ActiveWindow.Height = 15840

# ---------------------------------------------------------------------------
# 2. "stance t test" sheet loses its tabSelected flag (no longer the active
#    tab when the file is reopened) - achieved by activating another sheet
#    later in the script (alldata_1step becomes the selected tab).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 3. "alldata_1step" sheet: the stray AG:BE block (swing/stance scratch
#    tables) is removed - it was duplicated / relocated to "pythonout".
# ---------------------------------------------------------------------------
$wsAll1Step = $wb.Worksheets.Item("alldata_1step")
$wsAll1Step.Range("AG1:BE34").ClearContents()
$wsAll1Step.Range("F1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. "pythonout" sheet: append the swing_df and stance_df summary tables
#    (subject/condition/trial/value) that used to live in alldata_1step.
# ---------------------------------------------------------------------------
$wsPythonOut = $wb.Worksheets.Item("pythonout")

$swingRows = @(
    @("welk002","welkexo","trial01",2.8224860000000001),
    @("welk002","welkexo","trial02",2.7806150000000001),
    @("welk002","welkexo","trial03",2.5526460000000002),
    @("welk002","welkexo","trial04",2.8426930000000001),
    @("welk002","welknatural","trial01",2.670731),
    @("welk002","welknatural","trial02",3.1026220000000002),
    @("welk002","welknatural","trial03",2.6562060000000001),
    @("welk002","welknatural","trial04",2.774289),
    @("welk003","welkexo","trial01",3.190118),
    @("welk003","welkexo","trial02",3.1468889999999998),
    @("welk003","welkexo","trial03",3.2110289999999999),
    @("welk003","welkexo","trial04",3.5255000000000001),
    @("welk003","welknatural","trial01",3.3469850000000001),
    @("welk003","welknatural","trial02",3.0652469999999998),
    @("welk003","welknatural","trial03",3.0850390000000001),
    @("welk003","welknatural","trial04",3.1780020000000002)
)

$stanceRows = @(
    @("welk002","welkexo","trial01",6.6376109999999997),
    @("welk002","welkexo","trial02",5.9639850000000001),
    @("welk002","welkexo","trial03",6.6552610000000003),
    @("welk002","welkexo","trial04",5.8253329999999997),
    @("welk002","welknatural","trial01",7.4643920000000001),
    @("welk002","welknatural","trial02",7.6236660000000001),
    @("welk002","welknatural","trial03",6.9179539999999999),
    @("welk002","welknatural","trial04",6.8632220000000004),
    @("welk003","welkexo","trial01",6.9033720000000001),
    @("welk003","welkexo","trial02",6.647322),
    @("welk003","welkexo","trial03",6.3015169999999996),
    @("welk003","welkexo","trial04",6.898752),
    @("welk003","welknatural","trial01",7.3225720000000001),
    @("welk003","welknatural","trial02",7.7138470000000003),
    @("welk003","welknatural","trial03",7.5333889999999997),
    @("welk003","welknatural","trial04",7.8034150000000002)
)

# --- swing_df block: header label, column headers, then 16 data rows -------
$wsPythonOut.Cells.Item(41, 11).Value = "swing_df"

$wsPythonOut.Cells.Item(42, 12).Value = "subjectname"
$wsPythonOut.Cells.Item(42, 13).Value = "condname"
$wsPythonOut.Cells.Item(42, 14).Value = "trialname"
$wsPythonOut.Cells.Item(42, 15).Value = "metabolics_swing_avg_mean"

$r = 43
$i = 0
foreach ($row in $swingRows) {
    $wsPythonOut.Cells.Item($r, 11).Value = $i
    $wsPythonOut.Cells.Item($r, 12).Value = $row[0]
    $wsPythonOut.Cells.Item($r, 13).Value = $row[1]
    $wsPythonOut.Cells.Item($r, 14).Value = $row[2]
    $wsPythonOut.Cells.Item($r, 15).Value = $row[3]
    $r = $r + 1
    $i = $i + 1
}

# --- stance_df block: header label, column headers, then 16 data rows ------
$wsPythonOut.Cells.Item(59, 11).Value = "stance_df"

$wsPythonOut.Cells.Item(60, 12).Value = "subjectname"
$wsPythonOut.Cells.Item(60, 13).Value = "condname"
$wsPythonOut.Cells.Item(60, 14).Value = "trialname"
$wsPythonOut.Cells.Item(60, 15).Value = "metabolics_stance_avg_mean"

$r = 61
$i = 0
foreach ($row in $stanceRows) {
    $wsPythonOut.Cells.Item($r, 11).Value = $i
    $wsPythonOut.Cells.Item($r, 12).Value = $row[0]
    $wsPythonOut.Cells.Item($r, 13).Value = $row[1]
    $wsPythonOut.Cells.Item($r, 14).Value = $row[2]
    $wsPythonOut.Cells.Item($r, 15).Value = $row[3]
    $r = $r + 1
    $i = $i + 1
}

$wsPythonOut.Range("P43").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Make "alldata_1step" the active/selected sheet & set its view + final
#    selection as recorded in the saved workbook.
# ---------------------------------------------------------------------------
$wsAll1Step.Activate()
$wsAll1Step.Range("AG1:AQ16").Select() | Out-Null
